$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry updates one cell's displayed text. "ForceText" marks values that
# look numeric (e.g. "574.98", "0.510", "366.00") but must stay plain text so
# trailing/leading zeros and exact formatting survive, matching how the sheet
# stores prices/volumes as inline strings rather than numbers.
$updates = @(
    @{ Cell = 'D2'; Value = '66.615.85'; ForceText = $False },
    @{ Cell = 'D3'; Value = '3.070.89'; ForceText = $False },
    @{ Cell = 'E3'; Value = '  -1.39%  '; ForceText = $False },
    @{ Cell = 'E4'; Value = '  +0.05%  '; ForceText = $False },
    @{ Cell = 'D5'; Value = '574.98'; ForceText = $True },
    @{ Cell = 'E5'; Value = '  -0.57%  '; ForceText = $False },
    @{ Cell = 'D6'; Value = '168.92'; ForceText = $True },
    @{ Cell = 'E6'; Value = '  -1.70%  '; ForceText = $False },
    @{ Cell = 'E7'; Value = '  +0.11%  '; ForceText = $False },
    @{ Cell = 'D8'; Value = '3.067.06'; ForceText = $False },
    @{ Cell = 'E8'; Value = '  -1.40%  '; ForceText = $False },
    @{ Cell = 'D9'; Value = '0.510'; ForceText = $True },
    @{ Cell = 'E9'; Value = '  -2.15%  '; ForceText = $False },
    @{ Cell = 'D10'; Value = '6.42'; ForceText = $True },
    @{ Cell = 'E10'; Value = '  -0.23%  '; ForceText = $False },
    @{ Cell = 'E11'; Value = '  -1.85%  '; ForceText = $False },
    @{ Cell = 'D12'; Value = '0.467'; ForceText = $True },
    @{ Cell = 'E12'; Value = '  -3.51%  '; ForceText = $False },
    @{ Cell = 'D13'; Value = '0.0000238'; ForceText = $True },
    @{ Cell = 'E13'; Value = '  -2.79%  '; ForceText = $False },
    @{ Cell = 'D14'; Value = '35.57'; ForceText = $True },
    @{ Cell = 'E14'; Value = '  -4.12%  '; ForceText = $False },
    @{ Cell = 'E15'; Value = '  -1.82%  '; ForceText = $False },
    @{ Cell = 'D16'; Value = '3.583.84'; ForceText = $False },
    @{ Cell = 'E16'; Value = '  -1.26%  '; ForceText = $False },
    @{ Cell = 'D17'; Value = '66.607.33'; ForceText = $False },
    @{ Cell = 'E17'; Value = '  -0.28%  '; ForceText = $False },
    @{ Cell = 'D18'; Value = '16.90'; ForceText = $True },
    @{ Cell = 'E18'; Value = '  +3.67%  '; ForceText = $False },
    @{ Cell = 'D19'; Value = '6.95'; ForceText = $True },
    @{ Cell = 'E19'; Value = '  -3.24%  '; ForceText = $False },
    @{ Cell = 'D20'; Value = '3.077.01'; ForceText = $False },
    @{ Cell = 'E20'; Value = '  -1.17%  '; ForceText = $False },
    @{ Cell = 'D21'; Value = '486.86'; ForceText = $True },
    @{ Cell = 'E21'; Value = '  +2.28%  '; ForceText = $False },
    @{ Cell = 'D22'; Value = '7.68'; ForceText = $True },
    @{ Cell = 'E22'; Value = '  -2.51%  '; ForceText = $False },
    @{ Cell = 'D23'; Value = '0.685'; ForceText = $True },
    @{ Cell = 'E23'; Value = '  -3.98%  '; ForceText = $False },
    @{ Cell = 'D24'; Value = '82.57'; ForceText = $True },
    @{ Cell = 'E24'; Value = '  -1.84%  '; ForceText = $False },
    @{ Cell = 'D25'; Value = '12.65'; ForceText = $True },
    @{ Cell = 'E25'; Value = '  -4.85%  '; ForceText = $False },
    @{ Cell = 'D26'; Value = '2.20'; ForceText = $True },
    @{ Cell = 'E26'; Value = '  -3.60%  '; ForceText = $False },
    @{ Cell = 'D27'; Value = '10.11'; ForceText = $True },
    @{ Cell = 'E27'; Value = '  -0.74%  '; ForceText = $False },
    @{ Cell = 'E29'; Value = '  -1.32%  '; ForceText = $False },
    @{ Cell = 'E30'; Value = '  -4.80%  '; ForceText = $False },
    @{ Cell = 'D31'; Value = '2.60'; ForceText = $True },
    @{ Cell = 'E31'; Value = '  -2.75%  '; ForceText = $False },
    @{ Cell = 'D32'; Value = '27.49'; ForceText = $True },
    @{ Cell = 'E32'; Value = '  -3.78%  '; ForceText = $False },
    @{ Cell = 'E33'; Value = '  -3.86%  '; ForceText = $False },
    @{ Cell = 'D34'; Value = '0.0₃0905'; ForceText = $False },
    @{ Cell = 'E34'; Value = '  -3.88%  '; ForceText = $False },
    @{ Cell = 'E35'; Value = '  -0.01%  '; ForceText = $False },
    @{ Cell = 'D36'; Value = '5.58'; ForceText = $True },
    @{ Cell = 'E36'; Value = '  -4.89%  '; ForceText = $False },
    @{ Cell = 'D37'; Value = '0.946'; ForceText = $True },
    @{ Cell = 'E37'; Value = '  -2.67%  '; ForceText = $False },
    @{ Cell = 'D38'; Value = '46.92'; ForceText = $True },
    @{ Cell = 'E38'; Value = '  +0.03%  '; ForceText = $False },
    @{ Cell = 'E39'; Value = '  -0.71%  '; ForceText = $False },
    @{ Cell = 'E40'; Value = '  -5.05%  '; ForceText = $False },
    @{ Cell = 'E41'; Value = '  -3.51%  '; ForceText = $False },
    @{ Cell = 'D42'; Value = '8.27'; ForceText = $True },
    @{ Cell = 'E42'; Value = '  -4.97%  '; ForceText = $False },
    @{ Cell = 'D43'; Value = '2.760.35'; ForceText = $False },
    @{ Cell = 'E43'; Value = '  -2.55%  '; ForceText = $False },
    @{ Cell = 'E44'; Value = '  -2.63%  '; ForceText = $False },
    @{ Cell = 'B45'; Value = 'VeChain'; ForceText = $False },
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; ForceText = $False },
    @{ Cell = 'D45'; Value = '0.0344'; ForceText = $True },
    @{ Cell = 'E45'; Value = '  -3.34%  '; ForceText = $False },
    @{ Cell = 'B46'; Value = 'Monero'; ForceText = $False },
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; ForceText = $False },
    @{ Cell = 'D46'; Value = '135.35'; ForceText = $True },
    @{ Cell = 'E46'; Value = '  -0.08%  '; ForceText = $False },
    @{ Cell = 'D47'; Value = '366.00'; ForceText = $True },
    @{ Cell = 'E47'; Value = '  -5.63%  '; ForceText = $False },
    @{ Cell = 'D49'; Value = '24.60'; ForceText = $True },
    @{ Cell = 'E49'; Value = '  -0.88%  '; ForceText = $False },
    @{ Cell = 'E50'; Value = '  -2.15%  '; ForceText = $False },
    @{ Cell = 'E51'; Value = '  -2.09%  '; ForceText = $False }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Leading apostrophe forces Excel to treat the literal as text instead
        # of coercing it to a number (which would drop formatting like the
        # trailing zero in "574.98" -> 574.98 vs "46.92" staying "46.92").
        $rng.Value = "'" + $u.Value
        $rng.Style = "Normal"
    } else {
        $rng.Value = $u.Value
    }
}
